$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new ingredient row (row 49) with "water" / "Water" and zeroed nutrition values
$ws.Range("A49").Value = "water"
$ws.Range("B49").Value = "Water"
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0

# Update the view state to match the scrolled-down position after the edit
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B50").Select()
